# Entrega2 Planilha teste Local do resultado
#
# The "Resultado do teste" text that used to live in column C is moved to
# column D (its proper header column); column C is cleared out except for a
# leftover single-space placeholder in C3. Row 3, row 6 and row 10 also get
# their text content tweaked, row 6 no longer needs the tall wrapped row, and
# the view/column widths/selection are refreshed to match where the author
# left off editing (cell D10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move "Resultado do teste" values from column C to column D ----------

# Row 2: unchanged text, just relocated C2 -> D2
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "Falha SL:2275,00 IR=0,0"

# Row 3: relocated C3 -> D3, trailing space dropped from the moved text;
# C3 itself keeps a lone-space placeholder behind.
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = "Falha SL:1380,00"

# Row 4: unchanged text, just relocated C4 -> D4
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "Falha SL:1638,00  "

# Row 5: unchanged text, just relocated C5 -> D5
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "Falha SL:3066,00 IR=48,17 "

# Row 6: relocated C6 -> D6, trailing newline+space dropped from the moved
# text; row no longer needs to be tall/wrapped.
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "Falha SL: O valor do salario não pode ser menor que o salario minimo!!!!!Digite o valor do salario novamente:"
$ws.Range("C6").WrapText = $false
$ws.Rows.Item(6).RowHeight = 13.8

# Row 7: unchanged text, just relocated C7 -> D7
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = "Falha SL: O valor do salario não pode ser menor que o salario minimo!!!!!Digite o valor do salario novamente:  "

# Row 8: unchanged text, just relocated C8 -> D8
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = "Falha SL: O valor do salario não pode ser menor que o salario minimo!!!!!Digite o valor do salario novamente:"

# Row 9: unchanged text, just relocated C9 -> D9
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = "Falha SL: O valor do salario não pode ser menor que o salario minimo!!!!!Digite o valor do salario novamente:"

# Row 10: relocated C10 -> D10 with a brand new result string
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = "Falha= SL: 1112,00"

# --- Column widths: slight narrowing of A-D (matches saved state) --------

$ws.Columns.Item(1).ColumnWidth = 42.14625850340136
$ws.Columns.Item(2).ColumnWidth = 45.314625850340164
$ws.Columns.Item(3).ColumnWidth = 19.396258503401366
$ws.Columns.Item(4).ColumnWidth = 16.554421768707467

# --- View state: selection moved to D10, scrolled back to column A -------

$ws.Range("D10").Select()
